$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("M18").Value = -716

$ws.Range("H32").Value = 5068.136
$ws.Range("I32").Value = 9999.666999999999
$ws.Range("J32").Value = 4289.4736
$ws.Range("K32").Value = 9999.666999999999
$ws.Range("L32").Value = 4289.4736
$ws.Range("M32").Value = -9673.666999999999
$ws.Range("N32").Value = -4941.4736

$ws.Range("H87").Value = 68000
$ws.Range("J87").Value = 68000
$ws.Range("L87").Value = 68000
$ws.Range("N87").Value = -70496

$ws.Range("H90").Value = 68000
$ws.Range("J90").Value = 68000
$ws.Range("L90").Value = 204000
$ws.Range("N90").Value = -216480

$ws.Range("H101").Value = 2993.889
$ws.Range("I101").Value = 549.4
$ws.Range("J101").Value = 6049.5
$ws.Range("K101").Value = 1648.2
$ws.Range("L101").Value = 18148.5
$ws.Range("M101").Value = -26.19999999999982
$ws.Range("N101").Value = -21392.5

$ws.Range("H137").Value = 1388.6154
$ws.Range("I137").Value = 1245.2222
$ws.Range("J137").Value = 1711.25
$ws.Range("K137").Value = 3735.6666
$ws.Range("L137").Value = 5133.75
$ws.Range("M137").Value = -1185.6666
$ws.Range("N137").Value = -10233.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4942.3
$ws.Range("I2").Value = 4551.25
$ws.Range("K2").Value = 4551.25
$ws.Range("M2").Value = -4438.25

$ws.Range("H55").Value = 24770
$ws.Range("J55").Value = 34640
$ws.Range("L55").Value = 34640
$ws.Range("N55").Value = -35270

$ws.Range("H62").Value = 20000
$ws.Range("J62").Value = 20000
$ws.Range("L62").Value = 20000
$ws.Range("N62").Value = -21248

$ws.Range("H65").Value = 20000
$ws.Range("J65").Value = 20000
$ws.Range("L65").Value = 60000
$ws.Range("N65").Value = -66240

$ws.Range("H116").Value = 4942.3
$ws.Range("I116").Value = 4551.25
$ws.Range("K116").Value = 4551.25
$ws.Range("M116").Value = -2257.25

$ws.Range("H122").Value = 2075.375
$ws.Range("I122").Value = 1600.25
$ws.Range("J122").Value = 3500.75
$ws.Range("K122").Value = 4800.75
$ws.Range("L122").Value = 10502.25
$ws.Range("M122").Value = -2350.75
$ws.Range("N122").Value = -15402.25

$ws.Range("H132").Value = 6326.212
$ws.Range("I132").Value = 5411.1665
$ws.Range("K132").Value = 16233.4995
$ws.Range("M132").Value = -13703.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4942.3
$ws.Range("I3").Value = 4551.25
$ws.Range("K3").Value = 4551.25
$ws.Range("M3").Value = -4437.25

$ws.Range("H35").Value = 31750
$ws.Range("J35").Value = 31750
$ws.Range("L35").Value = 31750
$ws.Range("N35").Value = -32370

$ws.Range("H82").Value = 31064.25
$ws.Range("I82").Value = 7257
$ws.Range("J82").Value = 39000
$ws.Range("K82").Value = 7257
$ws.Range("L82").Value = 39000
$ws.Range("M82").Value = -6874
$ws.Range("N82").Value = -39766

$ws.Range("H85").Value = 31064.25
$ws.Range("I85").Value = 7257
$ws.Range("J85").Value = 39000
$ws.Range("K85").Value = 7257
$ws.Range("L85").Value = 39000
$ws.Range("M85").Value = -5931
$ws.Range("N85").Value = -41652

$ws.Range("H105").Value = 4000
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 4002
$ws.Range("I2").Value = 4002
$ws.Range("K2").Value = 4002
$ws.Range("M2").Value = -3889

$ws.Range("H7").Value = 510.76923
$ws.Range("I7").Value = 203
$ws.Range("K7").Value = 203
$ws.Range("M7").Value = -90

$ws.Range("H31").Value = 3794278.8
$ws.Range("I31").Value = 7581538.5
$ws.Range("J31").Value = 7018.8335
$ws.Range("K31").Value = 7581538.5
$ws.Range("L31").Value = 7018.8335
$ws.Range("M31").Value = -7581243.5
$ws.Range("N31").Value = -7608.8335

$ws.Range("H34").Value = 3794278.8
$ws.Range("I34").Value = 7581538.5
$ws.Range("J34").Value = 7018.8335
$ws.Range("K34").Value = 7581538.5
$ws.Range("L34").Value = 7018.8335
$ws.Range("M34").Value = -7581336.5
$ws.Range("N34").Value = -7422.8335

$ws.Range("H41").Value = 10627.143
$ws.Range("I41").Value = 7150
$ws.Range("J41").Value = 19320
$ws.Range("K41").Value = 7150
$ws.Range("L41").Value = 19320
$ws.Range("M41").Value = -6722
$ws.Range("N41").Value = -20176

$ws.Range("H58").Value = 9137.5
$ws.Range("J58").Value = 9371.429
$ws.Range("L58").Value = 9371.429
$ws.Range("N58").Value = -9777.429

$ws.Range("H60").Value = 24495
$ws.Range("I60").Value = 21242.5
$ws.Range("J60").Value = 31000
$ws.Range("K60").Value = 21242.5
$ws.Range("L60").Value = 31000
$ws.Range("M60").Value = -20731.5
$ws.Range("N60").Value = -32022

$ws.Range("H99").Value = 6840.727
$ws.Range("I99").Value = 6699.857
$ws.Range("K99").Value = 6699.857
$ws.Range("M99").Value = -5201.857

$ws.Range("H126").Value = 6840.727
$ws.Range("I126").Value = 6699.857
$ws.Range("K126").Value = 20099.571
$ws.Range("M126").Value = -17629.571

$ws.Range("H134").Value = 5280.4443
$ws.Range("I134").Value = 2574.9285
$ws.Range("J134").Value = 14749.75
$ws.Range("K134").Value = 7724.7855
$ws.Range("L134").Value = 44249.25
$ws.Range("M134").Value = -5189.7855
$ws.Range("N134").Value = -49319.25

$ws.Range("H136").Value = 9137.5
$ws.Range("J136").Value = 9371.429
$ws.Range("L136").Value = 28114.287
$ws.Range("N136").Value = -33214.287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 5172.6
$ws.Range("J2").Value = 10183.7
$ws.Range("L2").Value = 61102.2
$ws.Range("N2").Value = -61328.2

$ws.Range("H92").Value = 399.75
$ws.Range("J92").Value = 449.5
$ws.Range("L92").Value = 1348.5
$ws.Range("N92").Value = -3844.5

$ws.Range("H122").Value = 805.5
$ws.Range("I122").Value = 199
$ws.Range("J122").Value = 1007.6667
$ws.Range("K122").Value = 1791
$ws.Range("L122").Value = 9069.0003
$ws.Range("M122").Value = 659
$ws.Range("N122").Value = -13969.0003

$ws.Range("H132").Value = 33334600
$ws.Range("J132").Value = 1384.6154
$ws.Range("L132").Value = 12461.5386
$ws.Range("N132").Value = -17521.5386

$ws.Range("H137").Value = 22408.334
$ws.Range("I137").Value = 697.5
$ws.Range("J137").Value = 65830
$ws.Range("K137").Value = 2092.5
$ws.Range("L137").Value = 197490
$ws.Range("M137").Value = 3007.5
$ws.Range("N137").Value = -207690

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 20833.5
$ws.Range("I113").Value = 29875.25
$ws.Range("K113").Value = 29875.25
$ws.Range("M113").Value = -27705.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 7162.3335
$ws.Range("I93").Value = 2015.6666
$ws.Range("J93").Value = 9735.666999999999
$ws.Range("K93").Value = 2015.6666
$ws.Range("L93").Value = 9735.666999999999
$ws.Range("M93").Value = -767.6666
$ws.Range("N93").Value = -12231.667

$ws.Range("H132").Value = 8649.73
$ws.Range("I132").Value = 9800.368
$ws.Range("J132").Value = 5526.5713
$ws.Range("K132").Value = 29401.104
$ws.Range("L132").Value = 16579.7139
$ws.Range("M132").Value = -26871.104
$ws.Range("N132").Value = -21639.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 49999
$ws.Range("J120").Value = 49999
$ws.Range("L120").Value = 49999
$ws.Range("N120").Value = -59675

$ws.Range("H132").Value = 6240.9546
$ws.Range("J132").Value = 13333.333
$ws.Range("L132").Value = 39999.999
$ws.Range("N132").Value = -45059.999
